$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data change: C2 fraction value goes from 0.1 to 1E-3 (0.001) ---
$ws.Range("C2").Value = 0.001

# --- Column layout change: drop the old per-column bestFit widths on A:B
#     and give A:C a single, uniform custom width instead ---
$ws.Columns.Item(1).ColumnWidth = 19.33
$ws.Columns.Item(2).ColumnWidth = 19.33
$ws.Columns.Item(3).ColumnWidth = 19.33

# --- Selection change: active cell moves from C3 to B12 ---
$ws.Range("B12").Select()
